$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 was a duplicate "fortran-lang/fpm" contract entry (same org/repo/url
# info as row 4, differing only in the volume/date/contracting_entity
# columns). It was removed from the dataset, which also drops the
# now-orphaned "NumFocus (Fortran (2023))" shared string and shifts every
# later row up by one.
$ws.Rows.Item(5).Delete()

# The worksheet's hyperlinks (anchored to specific D-column cells) don't
# automatically re-anchor when the row shift happens, so re-create them
# pointing at their new (shifted-up-by-one) rows, preserving the same
# target URLs and the "Hyperlink" cell style.
$ws.Cells.Item(1, 1).Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Cells.Item(7, 4), "https://github.com/Lullabot/drupal9ci") | Out-Null
$ws.Cells.Item(7, 4).Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Cells.Item(25, 4), "https://github.com/sequoia-pgp/fast-forward") | Out-Null
$ws.Cells.Item(25, 4).Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Cells.Item(6, 4), "https://github.com/GStreamer/qt-gstreamer") | Out-Null
$ws.Cells.Item(6, 4).Style = "Hyperlink"

# Update the on-screen selection/scroll position to match the post-edit view.
$ws.Range("A4").Select()
